$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "in"
$ws2 = $wb.Worksheets.Item(2)   # "out"

# Fix the typo in the header: "conkat" -> "concat"
$ws1.Range("F1").Value = "concat"

# The aggregate row (Sum/Min/Max/Concat over the table) can no longer be relied
# upon when the table has no matching rows, so drop the computed formulas -
# leave D8 as an empty, still-formatted cell and remove B8/C8/F8 entirely.
[void]$ws1.Range("B8:F8").ClearContents()

# The "out" sheet held a cached copy of the old (typo'd) concatenation result;
# refresh it to match the corrected text now that the formula is gone.
$ws2.Range("D1").Value = "HelloOne WorldПривет twoмир!"

# Update the saved selections: "out" no longer is the active tab, "in" is -
# select "out"'s cell first so selecting on "in" afterwards leaves "in" active.
[void]$ws2.Range("D2").Select()
[void]$ws1.Range("F16").Select()
